# Apply cryptos list update (prices & 1h volume changes)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.916.34'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.49%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.254.54'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.20%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '583.50'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.74%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '183.87'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +3.96%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -1.23%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +3.55%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.69'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.51%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.415'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +1.64%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '3.820.78'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.04%  '
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.27%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.49'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +1.59%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '67.925.14'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +1.53%  '
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +2.67%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.254.43'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.03%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.85'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.72%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.59'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +1.45%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '381.56'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +3.27%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.67'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.86%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.07%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '71.30'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.76%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.73%  '
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.88%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.82'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.61%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +2.31%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.999'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.06%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.38%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.67'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.66%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.30'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +7.81%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '22.87'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +1.48%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +1.97%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +3.80%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '162.58'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -4.52%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.85'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.45%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.834'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -2.65%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '26.61'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -1.45%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.73'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +4.69%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +6.52%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.59'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.71%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '41.38'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +2.37%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '25.42'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +2.11%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '345.27'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +1.37%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.642.91'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -3.86%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0284'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +2.24%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.85%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.993'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +1.42%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '31.28'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +3.14%  '
